$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 17.17899368645011
    "D2" = 9.84383232753013
    "E2" = 16.03202598371829
    "F2" = 31.85767197334945
    "G2" = 3.645402074045877
    "I2" = 25.83248016663205
    "J2" = 11.43983118652399
    "L2" = 11.76141974793471
    "M2" = 17.06442421385967
    "O2" = 24.11367487535943
    "B3" = 16.72307153153398
    "D3" = 9.842262796629237
    "E3" = 15.93477218845757
    "F3" = 32.01102044269985
    "G3" = 3.647862164162222
    "I3" = 26.01299057684312
    "J3" = 11.39073110077423
    "L3" = 11.53563854483479
    "M3" = 16.8009781060472
    "O3" = 24.18808299262954
    "B4" = 16.43691650426402
    "D4" = 9.84198624809228
    "E4" = 15.87579671440022
    "F4" = 32.11498289667441
    "G4" = 3.649453320375299
    "I4" = 26.13025478076035
    "J4" = 11.36115144133652
    "L4" = 11.39501014264128
    "M4" = 16.63762204534032
    "O4" = 24.2410198246343
    "B5" = 16.31888802499297
    "D5" = 9.842047369690899
    "E5" = 15.85196515354093
    "F5" = 32.15980309975418
    "G5" = 3.650122076119619
    "I5" = 26.1796595167226
    "J5" = 11.34924693629256
    "L5" = 11.33726274248237
    "M5" = 16.57071519141923
    "O5" = 24.26440709029335
    "B6" = 16.29920827365999
    "D6" = 9.842068046641817
    "E6" = 15.84802055864174
    "F6" = 32.16739338745937
    "G6" = 3.650234353342264
    "I6" = 26.18796094904112
    "J6" = 11.34727941903941
    "L6" = 11.32764910345379
    "M6" = 16.55958684822972
    "O6" = 24.26839991156567
    "B7" = 16.43533027340115
    "D7" = 9.841986367269918
    "E7" = 15.87547447759882
    "F7" = 32.11557743364936
    "G7" = 3.64946225697829
    "I7" = 26.13091451422757
    "J7" = 11.36099027928283
    "L7" = 11.39423303958532
    "M7" = 16.63672099872906
    "O7" = 24.24132789524183
    "B8" = 17.02316631854674
    "D8" = 9.843149156971217
    "E8" = 15.99834901600979
    "F8" = 31.90850469538621
    "G8" = 3.646233614011429
    "I8" = 25.89338660892408
    "J8" = 11.42278760040132
    "L8" = 11.68402166269689
    "M8" = 16.97395797286022
    "O8" = 24.13782141011171
    "B9" = 18.12068573458077
    "D9" = 9.850836425281736
    "E9" = 16.24447085262033
    "F9" = 31.58072899571033
    "G9" = 3.640539158761563
    "I9" = 25.47853879469604
    "J9" = 11.54818384441733
    "L9" = 12.23383082159377
    "M9" = 17.61984647648524
    "O9" = 23.99270645041128
    "B10" = 18.88616893633917
    "D10" = 9.859719817664248
    "E10" = 16.42750297231205
    "F10" = 31.38829606131432
    "G10" = 3.636739483884431
    "I10" = 25.20470427896006
    "J10" = 11.6424798318651
    "L10" = 12.62320390273853
    "M10" = 18.08142469076117
    "O10" = 23.92179111457305
    "B11" = 19.22421009539929
    "D11" = 9.864450304936403
    "E11" = 16.51104488539823
    "F11" = 31.31139464312317
    "G11" = 3.635093403462218
    "I11" = 25.08683218033361
    "J11" = 11.68576080127717
    "L11" = 12.79653128469236
    "M11" = 18.2878852587092
    "O11" = 23.89736157066051
    "B12" = 19.35065605495111
    "D12" = 9.866339558322897
    "E12" = 16.54270295329761
    "F12" = 31.28381444455076
    "G12" = 3.634481857674093
    "I12" = 25.04315902862781
    "J12" = 11.70219800905394
    "L12" = 12.86156967558975
    "M12" = 18.36550507791934
    "O12" = 23.88924222684423
    "B13" = 19.32349464051223
    "D13" = 9.865928338350363
    "E13" = 16.53588409965122
    "F13" = 31.28968563805839
    "G13" = 3.634613041649743
    "I13" = 25.0525220242377
    "J13" = 11.69865596739865
    "L13" = 12.84758981185894
    "M13" = 18.34881415023344
    "O13" = 23.89094046254147
    "B14" = 19.23464473646822
    "D14" = 9.864603779514189
    "E14" = 16.51364906268819
    "F14" = 31.3090946634584
    "G14" = 3.635042855300509
    "I14" = 25.08321987718854
    "J14" = 11.68711219003173
    "L14" = 12.80189427305029
    "M14" = 18.29428267492245
    "O14" = 23.8966708787314
    "B15" = 19.18001529511059
    "D15" = 9.863805164691053
    "E15" = 16.50003187442011
    "F15" = 31.3211842193459
    "G15" = 3.635307662051316
    "I15" = 25.10214852606877
    "J15" = 11.68004726652446
    "L15" = 12.77382525644138
    "M15" = 18.26080569368538
    "O15" = 23.90032844639159
    "B16" = 18.86386361739981
    "D16" = 9.859424438383224
    "E16" = 16.42204743005279
    "F16" = 31.39353675314172
    "G16" = 3.636848712165926
    "I16" = 25.21254220331412
    "J16" = 11.63965840955179
    "L16" = 12.61179569272587
    "M16" = 18.06785636523744
    "O16" = 23.92354570217967
    "B17" = 18.66723386307352
    "D17" = 9.856912718974577
    "E17" = 16.37426556850805
    "F17" = 31.44065513420163
    "G17" = 3.637815160023186
    "I17" = 25.28197989504544
    "J17" = 11.61497432987946
    "L17" = 12.5113855338236
    "M17" = 17.94854683420857
    "O17" = 23.93979844995181
    "B18" = 18.55318464750045
    "D18" = 9.85553300704429
    "E18" = 16.34681008859752
    "F18" = 31.46875697039628
    "G18" = 3.638378795799737
    "I18" = 25.32254895440751
    "J18" = 11.60081343044169
    "L18" = 12.45327737016322
    "M18" = 17.87959567234967
    "O18" = 23.94988332122891
    "B19" = 18.51440903423517
    "D19" = 9.855077054195524
    "E19" = 16.3375193760305
    "F19" = 31.47844327890351
    "G19" = 3.63857096794773
    "I19" = 25.33639321974636
    "J19" = 11.59602533688631
    "L19" = 12.43354357984391
    "M19" = 17.85619553125511
    "O19" = 23.95342423564516
    "B20" = 18.68826481204642
    "D20" = 9.857173381111473
    "E20" = 16.37934932314419
    "F20" = 31.43553566354026
    "G20" = 3.637711477258581
    "I20" = 25.27452290101123
    "J20" = 11.61759824682025
    "L20" = 12.52211147255305
    "M20" = 17.96128185161433
    "O20" = 23.93799202649246
    "B21" = 19.26078526521742
    "D21" = 9.864990186811387
    "E21" = 16.52017956045982
    "F21" = 31.30335186150879
    "G21" = 3.634916289146327
    "I21" = 25.07417705617809
    "J21" = 11.69050164632867
    "L21" = 12.81533274694627
    "M21" = 18.31031559364713
    "O21" = 23.89495696399788
    "B22" = 19.62581137918558
    "D22" = 9.870669095988735
    "E22" = 16.61234357976426
    "F22" = 31.22594737455116
    "G22" = 3.633158160085003
    "I22" = 24.94884926226847
    "J22" = 11.73842230145364
    "L22" = 13.00346932439382
    "M22" = 18.53512715028513
    "O22" = 23.87342898495979
    "B23" = 19.43185749138867
    "D23" = 9.867586391592823
    "E23" = 16.56314845226962
    "F23" = 31.26643388395128
    "G23" = 3.63409024188205
    "I23" = 25.01522589370262
    "J23" = 11.71282364539021
    "L23" = 12.90339341903083
    "M23" = 18.41546129361615
    "O23" = 23.88431348117657
    "B24" = 18.67875984009478
    "D24" = 9.857055335365597
    "E24" = 16.37705090977556
    "F24" = 31.43784702203675
    "G24" = 3.63775832727395
    "I24" = 25.27789218818605
    "J24" = 11.61641187898128
    "L24" = 12.51726345859811
    "M24" = 17.95552546258332
    "O24" = 23.9388064029071
    "B25" = 17.83047176041756
    "D25" = 9.848184120584911
    "E25" = 16.17743204357376
    "F25" = 31.66095068658954
    "G25" = 3.642011913955803
    "I25" = 25.58532339007886
    "J25" = 11.51385077970331
    "L25" = 12.08742289714456
    "M25" = 17.61984647648524
    "O25" = 23.99270645041128
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

Write-Output "Updated $($values.Count) cells"
